$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the text "R40". The edit replaces its content
# with the text "1" while keeping it a text (shared-string) value rather
# than converting it to a number.
$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
